# Rename the existing sheet to the "CPU" variant and add a new "GPU" sheet
# with its own data + chart, mirroring the author's commit:
#   "Did some more testing and finished section on report about the Blur
#   kernel (minus the conclusion)."

$wb = $excel.ActiveWorkbook
$cpuWs = $wb.Worksheets.Item(1)
$cpuWs.Name = "Lena - Blur radius 1-15 - CPU"

# The chart on the CPU sheet still refers to the sheet by its old name;
# repoint its series at the renamed sheet (keeps all existing styling).
$cpuChart = $cpuWs.ChartObjects(1).Chart
$cpuSeries = $cpuChart.SeriesCollection(1)
$cpuSeries.Formula = "=SERIES('Lena - Blur radius 1-15 - CPU'!`$B`$1,'Lena - Blur radius 1-15 - CPU'!`$A`$2:`$A`$16,'Lena - Blur radius 1-15 - CPU'!`$B`$2:`$B`$16,1)"

# Fix up the CPU sheet's own view state (no longer the active tab; the
# author selected the full A:B columns before switching away).
$cpuWs.Range("A1:B1048576").Select()

# --- Add the GPU sheet -----------------------------------------------
$gpuWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cpuWs)
$gpuWs.Name = "Lena - Blur radius 1-15 - GPU"

$gpuWs.Range("A1").Value = "Blur Radius"
$gpuWs.Range("B1").Value = "Time"

$radii = 1..15
$times = @(
    0.0145992,
    0.046867699999999998,
    0.044389199999999997,
    0.068969199999999994,
    0.11994299999999999,
    0.158416,
    0.233708,
    0.29578100000000002,
    0.33590999999999999,
    0.39696999999999999,
    0.49583199999999999,
    0.57770999999999995,
    0.67206699999999997,
    0.75311700000000004,
    0.897818
)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $gpuWs.Cells.Item($row, 1).Value = $radii[$i]
    $gpuWs.Cells.Item($row, 2).Value = $times[$i]
}

$gpuWs.Columns.Item(1).ColumnWidth = 10.85546875

# --- Build the matching chart on the GPU sheet ------------------------
$gpuShape = $gpuWs.Shapes.AddChart2(-1, 4, 0, 0, 0, 0)
$gpuChart = $gpuShape.Chart
$gpuSeries = $gpuChart.SeriesCollection(1)
$gpuSeries.Formula = "=SERIES('Lena - Blur radius 1-15 - GPU'!`$B`$1,'Lena - Blur radius 1-15 - GPU'!`$A`$2:`$A`$16,'Lena - Blur radius 1-15 - GPU'!`$B`$2:`$B`$16,1)"

$gpuChart.HasTitle = $false
$gpuChart.HasLegend = $false

$catAxis = $gpuChart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Blur Radius"

$valAxis = $gpuChart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Time"

$gpuWs.Range("P13").Select()

Write-Host "Done"
